$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, copying the style of the existing
# header cell H1 ("IP") so the new headers match the other headers
# (bold, bordered, centered).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and column J (IF)
$ws.Range("I2").Value = 8
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I9").Value = 5

$ws.Range("J2").Value = 8
$ws.Range("J3").Value = 5
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 4
$ws.Range("J8").Value = 3
$ws.Range("J9").Value = 6
